# powerlaw fit.xlsx — add log(flux)/log(time) helper columns, a LINEST
# power-law fit (with/without the first two epochs), a second scatter
# chart with a linear trendline on the logged data, and the associated
# summary labels/values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Text labels first, in the exact order the original author must have
# entered them, so new shared-string indices line up (20..31).
# ---------------------------------------------------------------------
$ws.Range("K1").Value  = "log(flux)"
$ws.Range("L1").Value  = "log(time)"
$ws.Range("M15").Value = "rsq"
$ws.Range("M17").Value = "r"
$ws.Range("O2").Value  = "slope"
$ws.Range("P2").Value  = "y-int"
$ws.Range("N4").Value  = "errors"
$ws.Range("P16").Value = "Without Epoch 11"
$ws.Range("O18").Value = "Slope"
$ws.Range("O19").Value = "Error"
$ws.Range("O12").Value = "With all points"
$ws.Range("N14").Value = "error"
$ws.Range("N13").Value = "slope"
$ws.Range("P17").Value = "slope"

# ---------------------------------------------------------------------
# Per-row helper formulas: log10(flux) and log10(halfway-time)
# ---------------------------------------------------------------------
$ws.Range("K2").Formula  = "=LOG(G2:G13)"
$ws.Range("L2").Formula  = "=LOG(C2)"
$ws.Range("K3").Formula  = "=LOG(G3:G14)"
$ws.Range("L3").Formula  = "=LOG(C3)"
$ws.Range("K4").Formula  = "=LOG(G4:G15)"
$ws.Range("L4").Formula  = "=LOG(C4)"
$ws.Range("K5").Formula  = "=LOG(G5:G16)"
$ws.Range("L5").Formula  = "=LOG(C5)"
$ws.Range("K6").Formula  = "=LOG(G6:G17)"
$ws.Range("L6").Formula  = "=LOG(C6)"
$ws.Range("K7").Formula  = "=LOG(G7:G18)"
$ws.Range("L7").Formula  = "=LOG(C7)"
$ws.Range("K8").Formula  = "=LOG(G8:G19)"
$ws.Range("L8").Formula  = "=LOG(C8)"
$ws.Range("K9").Formula  = "=LOG(G9:G20)"
$ws.Range("L9").Formula  = "=LOG(C9)"
$ws.Range("K10").Formula = "=LOG(G10:G21)"
$ws.Range("L10").Formula = "=LOG(C10)"
$ws.Range("K11").Formula = "=LOG(G11:G22)"
$ws.Range("L11").Formula = "=LOG(C11)"
$ws.Range("K12").Formula = "=LOG(G12:G23)"
$ws.Range("L12").Formula = "=LOG(C12)"
$ws.Range("K13").Formula = "=LOG(G13:G24)"
$ws.Range("L13").Formula = "=LOG(C13)"

# ---------------------------------------------------------------------
# LINEST fit (without the first two epochs)
# ---------------------------------------------------------------------
$ws.Range("O3:P5").FormulaArray = "=LINEST(K4:K13,L4:L13,TRUE,TRUE)"
$ws.Range("O7").Formula = "=O3-O4"

# ---------------------------------------------------------------------
# rsq / r summary block
# ---------------------------------------------------------------------
$ws.Range("M16").Value = 0.8922
$ws.Range("M18").Formula = "=SQRT(M16)"

# ---------------------------------------------------------------------
# Condensed slope/error summary (typed copies of the LINEST results)
# ---------------------------------------------------------------------
$ws.Range("O13").Value = -1.9334270841585413
$ws.Range("O14").Value = 0.23763368633397353

# ---------------------------------------------------------------------
# "With all points" / "Without Epoch 11" comparison block
# ---------------------------------------------------------------------
$ws.Range("P18").Value = -2.0312621735276362
$ws.Range("P19").Value = 0.11134367972126213

# ---------------------------------------------------------------------
# Second chart: log(flux) vs log(time) scatter with a linear trendline
# ---------------------------------------------------------------------
$co2 = $ws.ChartObjects().Add(420, 110, 330, 260)
$chart2 = $co2.Chart
$chart2.ChartType = 74
$series2 = $chart2.SeriesCollection(1)
$series2.XValues = $ws.Range("L4:L13")
$series2.Values = $ws.Range("K4:K13")
$tl2 = $series2.Trendlines().Add()
$tl2.Type = -4132
$tl2.DisplayRSquared = $true
$tl2.DisplayEquation = $true

# ---------------------------------------------------------------------
# Selection, matching the saved view in the edited workbook
# ---------------------------------------------------------------------
$ws.Range("E13").Select()
